$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Atualiza notas dos alunos - preenche a coluna T6 (G) para cada aluno
$ws.Range("G2").Value = 1.1499999999999999
$ws.Range("G3").Value = 1.25
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 1.1000000000000001
$ws.Range("G6").Value = 1.25
$ws.Range("G6").NumberFormat = "0.00"

# Seleciona a célula G3, conforme estado final do arquivo
$ws.Range("G3").Select()
